$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "concepts"

$ws.Range("A1").Value = "concept"
$ws.Range("A2").Value = "gender"
$ws.Range("B1").Value = "nature"
$ws.Range("B2").Value = "entities"

$range = $ws.Range("A1:B2")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "concepts"
$tbl.HeaderRowRange.Font.Bold = $true

$ws.Columns.Item(1).ColumnWidth = 12.109375
$ws.Columns.Item(2).ColumnWidth = 11.21875

$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait
$ws.Range("C3").Select() | Out-Null
